# Updates the cryptos price list (columns D = Price, E = Volume(1h))
# for rows 2-51, matching the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values are plain text (e.g. "25.827.60" uses dots as thousand
# separators). Values that look like a single plain number (e.g. "215.86")
# would otherwise be auto-converted to a numeric cell by the Value setter,
# so those get a leading apostrophe to force them to stay text, exactly as
# Excel does when a user types a quote-prefixed numeric string.
function Set-TextValue($range, $value) {
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") "25.827.60"
$ws.Range("E2").Value = "  -0.42%  "
Set-TextValue $ws.Range("D3") "1.640.56"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "215.86"
$ws.Range("E5").Value = "  -0.11%  "
Set-TextValue $ws.Range("D6") "0.5063"
$ws.Range("E6").Value = "  -0.89%  "
Set-TextValue $ws.Range("D7") "1.004"
$ws.Range("E7").Value = "  +0.18%  "
Set-TextValue $ws.Range("D8") "0.2583"
$ws.Range("E8").Value = "  +0.33%  "
Set-TextValue $ws.Range("D9") "0.06438"
$ws.Range("E9").Value = "  +1.42%  "
Set-TextValue $ws.Range("D10") "20.42"
$ws.Range("E10").Value = "  +4.72%  "
Set-TextValue $ws.Range("D11") "0.07793"
$ws.Range("E11").Value = "  +0.27%  "
Set-TextValue $ws.Range("D12") "4.274"
$ws.Range("E12").Value = "  -0.26%  "
Set-TextValue $ws.Range("D13") "1.641.48"
$ws.Range("E13").Value = "  +0.21%  "
Set-TextValue $ws.Range("D14") "1.864.43"
$ws.Range("E14").Value = "  +0.29%  "
Set-TextValue $ws.Range("D15") "0.5624"
$ws.Range("E15").Value = "  +1.87%  "
Set-TextValue $ws.Range("D16") "0.0₅7659"
$ws.Range("E16").Value = "  +0.14%  "
Set-TextValue $ws.Range("D17") "63.44"
$ws.Range("E17").Value = "  -0.82%  "
Set-TextValue $ws.Range("D18") "25.839.28"
$ws.Range("E18").Value = "  -0.47%  "
Set-TextValue $ws.Range("D19") "1.003"
$ws.Range("E19").Value = "  +0.05%  "
Set-TextValue $ws.Range("D20") "4.387"
$ws.Range("E20").Value = "  -0.92%  "
Set-TextValue $ws.Range("D21") "193.17"
$ws.Range("E21").Value = "  -0.81%  "
Set-TextValue $ws.Range("D22") "9.944"
$ws.Range("E22").Value = "  +0.77%  "
Set-TextValue $ws.Range("D23") "6.150"
$ws.Range("E23").Value = "  +1.71%  "
Set-TextValue $ws.Range("D24") "1.003"
$ws.Range("E24").Value = "  +0.12%  "
Set-TextValue $ws.Range("D25") "1.799"
$ws.Range("E25").Value = "  -4.98%  "
Set-TextValue $ws.Range("D26") "141.21"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  -1.50%  "
Set-TextValue $ws.Range("D28") "6.833"
$ws.Range("E28").Value = "  +1.07%  "
Set-TextValue $ws.Range("D29") "15.61"
$ws.Range("E29").Value = "  +0.39%  "
Set-TextValue $ws.Range("D30") "1.246"
$ws.Range("E30").Value = "  +0.26%  "
Set-TextValue $ws.Range("D31") "0.04971"
$ws.Range("E31").Value = "  +2.00%  "
Set-TextValue $ws.Range("D32") "3.283"
$ws.Range("E32").Value = "  +1.26%  "
Set-TextValue $ws.Range("D33") "3.243"
$ws.Range("E33").Value = "  +1.50%  "
Set-TextValue $ws.Range("D34") "1.571"
$ws.Range("E34").Value = "  +1.89%  "
Set-TextValue $ws.Range("D35") "2.387"
$ws.Range("E35").Value = "  +0.68%  "
Set-TextValue $ws.Range("D36") "0.9055"
$ws.Range("E36").Value = "  +0.77%  "
Set-TextValue $ws.Range("D37") "2.570"
$ws.Range("E37").Value = "  +1.15%  "
Set-TextValue $ws.Range("D38") "0.5569"
$ws.Range("E38").Value = "  +0.95%  "
Set-TextValue $ws.Range("D39") "1.134.12"
$ws.Range("E39").Value = "  +1.33%  "
Set-TextValue $ws.Range("D40") "0.01571"
$ws.Range("E40").Value = "  +0.83%  "
Set-TextValue $ws.Range("D41") "0.9968"
$ws.Range("E41").Value = "  -0.48%  "
Set-TextValue $ws.Range("D42") "5.489"
$ws.Range("E42").Value = "  -1.84%  "
Set-TextValue $ws.Range("D43") "0.8026"
$ws.Range("E43").Value = "  +0.63%  "
Set-TextValue $ws.Range("D44") "98.99"
Set-TextValue $ws.Range("D45") "1.774.39"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("E46").Value = "  -6.50%  "
Set-TextValue $ws.Range("D47") "55.65"
$ws.Range("E47").Value = "  +1.58%  "
Set-TextValue $ws.Range("D48") "0.4270"
$ws.Range("E48").Value = "  -4.00%  "
Set-TextValue $ws.Range("D49") "7.760"
$ws.Range("E49").Value = "  +2.66%  "
Set-TextValue $ws.Range("D51") "0.9988"
$ws.Range("E51").Value = "  -0.39%  "
